$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.243.47"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.333.20"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.89"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.26"
$ws.Range("E6").Value = "  +4.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.360.94"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("E10").Value = "  +7.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  +6.23%  "
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.79"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.754.49"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.063.41"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.334.98"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.53"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.17"
$ws.Range("E21").Value = "  +5.05%  "
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.16"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  +7.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.85"
$ws.Range("E27").Value = "  +4.59%  "
$ws.Range("E28").Value = "  +9.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.71"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0748"
$ws.Range("E30").Value = "  +5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.69"
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.38"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.26"
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  +5.14%  "
$ws.Range("E39").Value = "  +8.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.90"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +5.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.12"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "279.63"
$ws.Range("E44").Value = "  +10.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.14"
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0506"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("E49").Value = "  +5.18%  "
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +10.98%  "
